$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of data (row 6)
$ws.Range("A6").Value = "John "
$ws.Range("B6").Value = "Dirty dishes"
$ws.Range("C6").Value = "wand"
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = "Have a nice day"

# Update selection to match new active cell
$ws.Range("I7").Select()
